$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.00005788376926308966
$ws.Range("E2").Value = 0.00005788376926308966

# Row 3
$ws.Range("D3").Value = 0.9691230853224901
$ws.Range("E3").Value = 0.9691230853224901

# Row 4
$ws.Range("D4").Value = 0.003988949099698905
$ws.Range("E4").Value = 0.003988949099698905

# Row 5
$ws.Range("D5").Value = 0.0001756157487735972
$ws.Range("E5").Value = 0.0001756157487735972

# Row 6
$ws.Range("D6").Value = 0.09952853282508131
$ws.Range("E6").Value = 0.09952853282508131

# Row 7
$ws.Range("D7").Value = 0.9999999863758322
$ws.Range("E7").Value = 0.00000001362416779393527

# Row 8
$ws.Range("D8").Value = 0.5643412826926471
$ws.Range("E8").Value = 0.4356587173073529

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.3270565767246524
$ws.Range("E9").Value = 0.6729434232753475

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.1436111315283776
$ws.Range("E10").Value = 0.8563888684716223

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.005866871156122744
$ws.Range("E11").Value = 0.9941331288438773
$ws.Range("F11").Value = 1.23556113243103
$ws.Range("G11").Value = 0.6

# Row 12
$ws.Range("D12").Value = 0.9956673612355355
$ws.Range("E12").Value = 0.9956673612355355

# Row 13
$ws.Range("D13").Value = 0.01129548860672427
$ws.Range("E13").Value = 0.01129548860672427

# Row 14
$ws.Range("D14").Value = 0.0004879512756157871
$ws.Range("E14").Value = 0.0004879512756157871

# Row 15
$ws.Range("D15").Value = 0.00000008011800586253572
$ws.Range("E15").Value = 0.00000008011800586253572

# Row 16
$ws.Range("D16").Value = 0.0217715204990546
$ws.Range("E16").Value = 0.0217715204990546

# Row 17
$ws.Range("D17").Value = 0.9999999999179285
$ws.Range("E17").Value = 0.0000000000820714607385753

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = 0.1413514940519286
$ws.Range("E18").Value = 0.8586485059480714

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = 0.3369531335800979
$ws.Range("E19").Value = 0.6630468664199021

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.02090072124396112
$ws.Range("E20").Value = 0.9790992787560389

# Row 21
$ws.Range("D21").Value = 0.9611249261613846
$ws.Range("E21").Value = 0.03887507383861544
$ws.Range("F21").Value = 1.242737889289856
$ws.Range("G21").Value = 0.6
